$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 431, shifting rows 431:462 down to 432:463
$ws.Rows.Item(431).Insert()

# Fill in the new row 431 with its values
$ws.Cells.Item(431, 1).Value = 5
$ws.Cells.Item(431, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(431, 3).Value = "Maule"
$ws.Cells.Item(431, 4).Value = 44783
$ws.Cells.Item(431, 5).Value = 7
$ws.Cells.Item(431, 6).Value = 100112043
$ws.Cells.Item(431, 7).Value = "Pepino ensalada"
$ws.Cells.Item(431, 8).Value = "Sin especificar"
$ws.Cells.Item(431, 9).Value = "Primera"
$ws.Cells.Item(431, 10).Value = 300
$ws.Cells.Item(431, 11).Value = 20000
$ws.Cells.Item(431, 12).Value = 20000
$ws.Cells.Item(431, 13).Value = 20000
$ws.Cells.Item(431, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(431, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(431, 16).Value = 333
$ws.Cells.Item(431, 17).Value = 60
$ws.Cells.Item(431, 18).Value = "Hortaliza"
